# "Saldo" export sheet ("Export"): the block of rows 5-16 (1-based Excel
# rows; row 1 is the header "Conta"/"Nome"/"Saldo") lists accounts together
# with their balance. The account 005437764 / EVA (balance 100000), which
# used to be the first of that block (row 5), needs to move down so it
# becomes the last row of the block (row 16) - i.e. the 11 rows that used to
# follow it (PEDRO, LAGO, MARCO, THOMAS, PAULO, GILTON, LUCAS, BIANCA,
# RICARDO, RAFAEL, FERNANDO) each shift up by one row and EVA's row ends up
# right before GUSTAVO's row (row 17, unchanged).
#
# We reproduce that reordering with plain Copy / PasteSpecial(xlPasteValues)
# operations (so cell formatting/styles are left exactly as Excel leaves
# them for an ordinary value copy - no new number-format styles get
# introduced), rather than relying on Cut/Insert of whole rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

$firstRow = 5
$lastRow = 16
$scratchRow = 10000

# 1. Stash EVA's row (currently row 5) out of the way.
$ws.Range("A$firstRow`:C$firstRow").Copy()
$ws.Range("A$scratchRow`:C$scratchRow").PasteSpecial($xlPasteValues)

# 2. Shift every row below it (6..16) up by one, in order.
for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $ws.Range("A$r`:C$r").Copy()
    $ws.Range("A$($r - 1)`:C$($r - 1)").PasteSpecial($xlPasteValues)
}

# 3. Drop EVA's stashed row into the now-vacant last slot of the block.
$ws.Range("A$scratchRow`:C$scratchRow").Copy()
$ws.Range("A$lastRow`:C$lastRow").PasteSpecial($xlPasteValues)

# 4. Clean up the scratch area and the marching-ants clipboard marquee.
$ws.Range("A$scratchRow`:C$scratchRow").Clear()
$excel.CutCopyMode = 0
